# Apply updated cryptocurrency price/volume data to cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.180.14'
$ws.Range("E2").Value = '  +1.82%  '

$ws.Range("D3").Value = '2.019.80'
$ws.Range("E3").Value = '  +3.32%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '''246.96'
$ws.Range("E5").Value = '  +1.33%  '

$ws.Range("E6").Value = '  +0.42%  '

$ws.Range("D7").Value = '''60.53'
$ws.Range("E7").Value = '  +0.65%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +3.77%  '

$ws.Range("D10").Value = '''0.0813'
$ws.Range("E10").Value = '  +2.81%  '

$ws.Range("E11").Value = '  +1.74%  '

$ws.Range("E12").Value = '  +6.07%  '

$ws.Range("D13").Value = '2.317.64'
$ws.Range("E13").Value = '  +3.48%  '

$ws.Range("E14").Value = '  +3.04%  '

$ws.Range("D15").Value = '''21.93'
$ws.Range("E15").Value = '  +1.73%  '

$ws.Range("E16").Value = '  +3.52%  '

$ws.Range("D17").Value = '2.019.47'
$ws.Range("E17").Value = '  +3.12%  '

$ws.Range("D18").Value = '37.119.47'
$ws.Range("E18").Value = '  +1.71%  '

$ws.Range("D19").Value = '''70.32'
$ws.Range("E19").Value = '  +1.59%  '

$ws.Range("D20").Value = '0.0₃0867'
$ws.Range("E20").Value = '  +1.72%  '

$ws.Range("D21").Value = '''5.22'
$ws.Range("E21").Value = '  +2.95%  '

$ws.Range("D22").Value = '''230.51'
$ws.Range("E22").Value = '  +0.51%  '

$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("D24").Value = '''2.58'
$ws.Range("E24").Value = '  +5.48%  '

$ws.Range("D25").Value = '''2.35'
$ws.Range("E25").Value = '  -0.46%  '

$ws.Range("D26").Value = '''9.39'
$ws.Range("E26").Value = '  +2.43%  '

$ws.Range("D27").Value = '''163.12'
$ws.Range("E27").Value = '  +1.88%  '

$ws.Range("E28").Value = '  -3.65%  '

$ws.Range("D29").Value = '''19.75'
$ws.Range("E29").Value = '  +2.50%  '

$ws.Range("D30").Value = '''1.39'
$ws.Range("E30").Value = '  +5.85%  '

$ws.Range("E31").Value = '  +1.02%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.0674'
$ws.Range("E32").Value = '  +10.13%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''4.77'
$ws.Range("E33").Value = '  +0.14%  '

$ws.Range("D34").Value = '''2.50'
$ws.Range("E34").Value = '  +10.16%  '

$ws.Range("D35").Value = '''4.45'
$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").Value = '''3.62'
$ws.Range("E36").Value = '  +5.55%  '

$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("D38").Value = '''1.81'
$ws.Range("E38").Value = '  +2.01%  '

$ws.Range("D39").Value = '''5.36'
$ws.Range("E39").Value = '  -1.59%  '

$ws.Range("E40").Value = '  +2.90%  '

$ws.Range("D41").Value = '''0.0977'
$ws.Range("E41").Value = '  +0.99%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.0215'
$ws.Range("E42").Value = '  +2.19%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '''16.89'
$ws.Range("E43").Value = '  +6.91%  '

$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '''1.18'
$ws.Range("E44").Value = '  +1.36%  '

$ws.Range("D45").Value = '''91.44'
$ws.Range("E45").Value = '  +3.07%  '

$ws.Range("D46").Value = '1.377.68'
$ws.Range("E46").Value = '  +0.98%  '

$ws.Range("E47").Value = '  +2.44%  '

$ws.Range("D48").Value = '''7.44'
$ws.Range("E48").Value = '  +4.28%  '

$ws.Range("E49").Value = '  +15.52%  '

$ws.Range("D50").Value = '''2.88'
$ws.Range("E50").Value = '  +1.71%  '

$ws.Range("D51").Value = '''46.15'
$ws.Range("E51").Value = '  +2.00%  '
